# "C program test files" — refresh the "Used Instructions" workbook:
#   - rename/extend the query-defined name to cover the larger imported range
#   - repopulate the imported instruction list (Sheet2, the data sheet) with
#     the new, longer set of MIPS mnemonics actually used by the C program
#   - a few cosmetic knock-on changes (selection, column width, page setup)

$wb = $excel.ActiveWorkbook

# --- Defined name: test_c_assembly -> test_c_assembly_1, wider range -----
try {
    $nm = $wb.Names.Item(1)
    $nm.Name = "test_c_assembly_1"
    $nm.RefersTo = "=Sheet2!`$A`$1:`$A`$145"
} catch {
    Write-Host "defined name update failed:" $_
}

# --- Data sheet (displayed as "Sheet2") -----------------------------------
$ws = $wb.Worksheets.Item("Sheet2")

$instructions = @(
    "addiu", "sw", "move", "li", "lw", "nop",
    "addu", "subu", "mult", "mflo", "sll", "bnez",
    "div", "break", "srl", "sra", "j", "bgtz", "slt", "jr"
)

for ($i = 0; $i -lt $instructions.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $instructions[$i]
}

# widen column A for the longer mnemonics, and relocate the selection
$ws.Columns.Item(1).ColumnWidth = 28.28515625
[void]$ws.Range("B15").Select()

# page setup now prints on A4/Letter-ish sheet (size id 9) in portrait
try {
    $ps = $ws.PageSetup
    $ps.PaperSize = 9
    $ps.Orientation = 1
} catch {
    Write-Host "page setup update failed:" $_
}

# --- Best-effort cosmetic/project metadata (may be no-ops on this host) --
try { $wb.CodeName = "ThisWorkbook" } catch {}
try { $ws.CodeName = "Sheet1" } catch {}
try { $wb.Worksheets.Item("Sheet1").CodeName = "Sheet2" } catch {}

try {
    $conn = $wb.Connections.Item(1)
    $tc = $conn.TextConnection
    $tc.TextFileCodePage = 850
    $tc.TextFilePlatform = 850
} catch {}

try {
    $qt = $ws.QueryTables.Item(1)
    $qt.Name = "test_c_assembly_1"
} catch {}
